# fix : 강화 2->3
# Normalize the Unit_Modeling (column K) values so every row within a
# Type/Step group uses the same model string, and fix the Mage group which
# was incorrectly pointing at the Archer model.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column K = Unit_Modeling. Rows 2-4 = 초보자(Base0), 5-7 = 전사(Human0),
# 8-10 = 궁수(Archer0), 11-13 = 마법사(Mage0).
$ws.Range("K2").Value = "Base0"
$ws.Range("K3").Value = "Base0"
$ws.Range("K4").Value = "Base0"

$ws.Range("K5").Value = "Human0"
$ws.Range("K6").Value = "Human0"
$ws.Range("K7").Value = "Human0"

$ws.Range("K8").Value = "Archer0"
$ws.Range("K9").Value = "Archer0"
$ws.Range("K10").Value = "Archer0"

$ws.Range("K11").Value = "Mage0"
$ws.Range("K12").Value = "Mage0"
$ws.Range("K13").Value = "Mage0"

# Cosmetic: the active cell selection moved from F16 to I16.
$ws.Range("I16").Select() | Out-Null
